$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row=2; E="3"; G="50.862619"; H="152.587857"; I="0.6466984659960481"; J="0.646698465996048"; K="3"; M="343.9479473333333"; N="1031.843842"; O="0.9666099193889262"; P="0.966609919388926"; Q="17494.0934010474"; R="157446.8406094266"; S="0.6251051520853822"; T="0.6251051520853821" },
    @{ Row=3; E="3"; G="50.862619"; H="152.587857"; I="0.6466984659960481"; J="0.646698465996048"; K="3"; M="5.970184"; N="17.910552"; O="0.01677823379880302"; P="0.01677823379880302"; Q="303.659194151896"; R="2732.932747367064"; S="0.01085045805980896"; T="0.01085045805980896" },
    @{ Row=4; E="3"; G="50.862619"; H="152.587857"; I="0.6466984659960481"; J="0.646698465996048"; K="3"; M="2.671367333333334"; N="8.014102000000001"; O="0.007507444608265281"; P="0.00750744460826528"; Q="135.8727388843793"; R="1222.854649959414"; S="0.00485505291171546"; T="0.004855052911715458" },
    @{ Row=5; E="3"; G="50.862619"; H="152.587857"; I="0.6466984659960481"; J="0.646698465996048"; K="3"; M="3.239611333333334"; N="9.718834000000001"; O="0.009104402204005551"; P="0.00910440220400555"; Q="164.7751169554153"; R="1482.976052598738"; S="0.005887802939141429"; T="0.005887802939141427" },
    @{ Row=6; E="3"; G="4.214243"; H="12.642729"; I="0.05358246462759977"; J="0.05358246462759976"; K="3"; M="343.9479473333333"; N="1031.843842"; O="0.9666099193889262"; P="0.966609919388926"; Q="1449.480229413869"; R="13045.32206472482"; S="0.0517933418143442"; T="0.05179334181434418" },
    @{ Row=7; E="3"; G="4.214243"; H="12.642729"; I="0.05358246462759977"; J="0.05358246462759976"; K="3"; M="5.970184"; N="17.910552"; O="0.01677823379880302"; P="0.01677823379880302"; Q="25.159806130712"; R="226.438255176408"; S="0.0008990191190379619"; T="0.0008990191190379618" },
    @{ Row=8; E="3"; G="4.214243"; H="12.642729"; I="0.05358246462759977"; J="0.05358246462759976"; K="3"; M="2.671367333333334"; N="8.014102000000001"; O="0.007507444608265281"; P="0.00750744460826528"; Q="11.25779108492867"; R="101.320119764358"; S="0.000402267385166039"; T="0.0004022673851660389" },
    @{ Row=9; E="3"; G="4.214243"; H="12.642729"; I="0.05358246462759977"; J="0.05358246462759976"; K="3"; M="3.239611333333334"; N="9.718834000000001"; O="0.009104402204005551"; P="0.00910440220400555"; Q="13.65250938422067"; R="122.872584457986"; S="0.0004878363090515688"; T="0.0004878363090515686" },
    @{ Row=10; E="3"; G="23.01971966666666"; H="69.05915899999999"; I="0.2926868039589623"; J="0.2926868039589623"; K="3"; M="343.9479473333333"; N="1031.843842"; O="0.9666099193889262"; P="0.966609919388926"; Q="7917.585327538764"; R="71258.26794784887"; S="0.282913967980975"; T="0.282913967980975" },
    @{ Row=11; E="3"; G="23.01971966666666"; H="69.05915899999999"; I="0.2926868039589623"; J="0.2926868039589623"; K="3"; M="5.970184"; N="17.910552"; O="0.01677823379880302"; P="0.01677823379880302"; Q="137.4319620384186"; R="1236.887658345768"; S="0.004910767626647897"; T="0.004910767626647896" },
    @{ Row=12; E="3"; G="23.01971966666666"; H="69.05915899999999"; I="0.2926868039589623"; J="0.2926868039589623"; K="3"; M="2.671367333333334"; N="8.014102000000001"; O="0.007507444608265281"; P="0.00750744460826528"; Q="61.49412714002423"; R="553.447144260218"; S="0.002197329968292109"; T="0.002197329968292108" },
    @{ Row=13; E="3"; G="23.01971966666666"; H="69.05915899999999"; I="0.2926868039589623"; J="0.2926868039589623"; K="3"; M="3.239611333333334"; N="9.718834000000001"; O="0.009104402204005551"; P="0.00910440220400555"; Q="74.57494472228956"; R="671.174502500606"; S="0.002664738383047318"; T="0.002664738383047316" },
    @{ Row=14; E="3"; G="0.5530853333333333"; H="1.659256"; I="0.007032265417389923"; J="0.007032265417389922"; K="3"; M="343.9479473333333"; N="1031.843842"; O="0.9666099193889262"; P="0.966609919388926"; Q="190.2325651001725"; R="1712.093085901552"; S="0.006797457508224807"; T="0.006797457508224805" },
    @{ Row=15; E="3"; G="0.5530853333333333"; H="1.659256"; I="0.007032265417389923"; J="0.007032265417389922"; K="3"; M="5.970184"; N="17.910552"; O="0.01677823379880302"; P="0.01677823379880302"; Q="3.302021207701333"; R="29.718190869312"; S="0.0001179889933082053"; T="0.0001179889933082052" },
    @{ Row=16; E="3"; G="0.5530853333333333"; H="1.659256"; I="0.007032265417389923"; J="0.007032265417389922"; K="3"; M="2.671367333333334"; N="8.014102000000001"; O="0.007507444608265281"; P="0.00750744460826528"; Q="1.477494092012445"; R="13.297446828112"; S="5.279434309167438E-05"; T="5.279434309167436E-05" },
    @{ Row=17; E="3"; G="0.5530853333333333"; H="1.659256"; I="0.007032265417389923"; J="0.007032265417389922"; K="3"; M="3.239611333333334"; N="9.718834000000001"; O="0.009104402204005551"; P="0.00910440220400555"; Q="1.791781514167111"; R="16.126033627504"; S="6.402457276523683E-05"; T="6.402457276523682E-05" }
)

$colIndex = @{ E=5; G=7; H=8; I=9; J=10; K=11; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20 }

foreach ($r in $rowData) {
    $rowNum = $r.Row
    foreach ($col in @('E','G','H','I','J','K','M','N','O','P','Q','R','S','T')) {
        $ws.Cells.Item($rowNum, $colIndex[$col]).Value = [double]$r[$col]
    }
}

Write-Output "Applied Cd38-Pecam1 NATMI edge updates"
